# Update the "想去人数" (F column) figures that were refreshed in the
# upstream data source. The same underlying events appear on both the
# "展览" sheet and the combined "全部类型" sheet, so both need updating.

$wb = $excel.ActiveWorkbook

$exhibitionUpdates = @{
    2  = 45
    3  = 777
    5  = 57
    6  = 59
    7  = 269
    8  = 3809
    10 = 4496
    12 = 1125
    13 = 62
}

$allTypesUpdates = @{
    2  = 45
    3  = 777
    5  = 57
    6  = 59
    8  = 269
    9  = 3809
    11 = 4496
    13 = 1125
    14 = 62
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
